# Ajout image assemblage + conversion
# Rename the three assembly ID codes from the old "STA000x" naming scheme
# to the new "STA0x00" naming scheme used after converting the assembly
# images.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "ST_A0100"
$ws.Range("G9").Value = "ST_A0400"
$ws.Range("G16").Value = "ST_A0500"

# Restore the view/selection state recorded in the saved workbook.
$ws.Range("G17").Select()
